$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-10 from 45212 to 45221
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45221
}
